$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Individual cell updates: Row(R), Column(C), Value(V), Force=text-forced via NumberFormat "@"
# (needed when the new value would otherwise be auto-parsed as a number by Excel,
# since the source cells are plain inline-string text, e.g. "0.9981", "240.43", etc.)
$edits = @(
    @{R=2; C=4; V='29.360.38'; Force=$false},
    @{R=2; C=5; V='  -0.20%  '; Force=$false},
    @{R=3; C=4; V='1.847.35'; Force=$false},
    @{R=3; C=5; V='  -0.14%  '; Force=$false},
    @{R=4; C=4; V='0.9981'; Force=$true},
    @{R=5; C=4; V='240.43'; Force=$true},
    @{R=5; C=5; V='  -0.15%  '; Force=$false},
    @{R=6; C=4; V='0.6272'; Force=$true},
    @{R=6; C=5; V='  -0.55%  '; Force=$false},
    @{R=7; C=4; V='0.9990'; Force=$true},
    @{R=7; C=5; V='  -0.16%  '; Force=$false},
    @{R=8; C=4; V='0.07622'; Force=$true},
    @{R=8; C=5; V='  -1.12%  '; Force=$false},
    @{R=10; C=5; V='  +1.11%  '; Force=$false},
    @{R=11; C=4; V='0.07734'; Force=$true},
    @{R=12; C=4; V='5.027'; Force=$true},
    @{R=12; C=5; V='  +0.04%  '; Force=$false},
    @{R=13; C=5; V='  -0.14%  '; Force=$false},
    @{R=14; C=4; V='0.00001051'; Force=$true},
    @{R=14; C=5; V='  -3.27%  '; Force=$false},
    @{R=15; C=4; V='82.98'; Force=$true},
    @{R=15; C=5; V='  -0.91%  '; Force=$false},
    @{R=16; C=4; V='6.145'; Force=$true},
    @{R=16; C=5; V='  -0.16%  '; Force=$false},
    @{R=17; C=4; V='29.381.95'; Force=$false},
    @{R=17; C=5; V='  -0.20%  '; Force=$false},
    @{R=18; C=4; V='227.93'; Force=$true},
    @{R=18; C=5; V='  -0.66%  '; Force=$false},
    @{R=19; C=5; V='  -1.04%  '; Force=$false},
    @{R=20; C=4; V='0.9985'; Force=$true},
    @{R=21; C=4; V='7.472'; Force=$true},
    @{R=21; C=5; V='  +0.21%  '; Force=$false},
    @{R=22; C=4; V='0.9990'; Force=$true},
    @{R=22; C=5; V='  -0.19%  '; Force=$false},
    @{R=23; C=4; V='158.85'; Force=$true},
    @{R=23; C=5; V='  +0.92%  '; Force=$false},
    @{R=24; C=5; V='  -0.28%  '; Force=$false},
    @{R=25; C=4; V='8.428'; Force=$true},
    @{R=25; C=5; V='  +0.75%  '; Force=$false},
    @{R=26; C=4; V='17.65'; Force=$true},
    @{R=26; C=5; V='  -0.18%  '; Force=$false},
    @{R=27; C=4; V='1.409'; Force=$true},
    @{R=27; C=5; V='  +7.48%  '; Force=$false},
    @{R=28; C=4; V='1.460'; Force=$true},
    @{R=28; C=5; V='  -0.60%  '; Force=$false},
    @{R=29; C=4; V='0.05596'; Force=$true},
    @{R=29; C=5; V='  -2.45%  '; Force=$false},
    @{R=30; C=4; V='4.106'; Force=$true},
    @{R=30; C=5; V='  -0.14%  '; Force=$false},
    @{R=31; C=5; V='  +0.25%  '; Force=$false},
    @{R=32; C=5; V='  +0.17%  '; Force=$false},
    @{R=33; C=4; V='1.834'; Force=$true},
    @{R=33; C=5; V='  -1.03%  '; Force=$false},
    @{R=34; C=4; V='0.6961'; Force=$true},
    @{R=34; C=5; V='  -1.86%  '; Force=$false},
    @{R=35; C=4; V='2.585'; Force=$true},
    @{R=35; C=5; V='  -0.12%  '; Force=$false},
    @{R=36; C=5; V='  +0.22%  '; Force=$false},
    @{R=37; C=4; V='1.228.23'; Force=$false},
    @{R=37; C=5; V='  -0.13%  '; Force=$false},
    @{R=38; C=5; V='  -2.09%  '; Force=$false},
    @{R=39; C=4; V='6.357'; Force=$true},
    @{R=39; C=5; V='  -1.94%  '; Force=$false},
    @{R=40; C=4; V='0.9029'; Force=$true},
    @{R=40; C=5; V='  -1.11%  '; Force=$false},
    @{R=41; C=4; V='0.9989'; Force=$true},
    @{R=41; C=5; V='  -0.18%  '; Force=$false},
    @{R=42; C=4; V='101.24'; Force=$true},
    @{R=42; C=5; V='  -0.36%  '; Force=$false},
    @{R=43; C=4; V='65.56'; Force=$true},
    @{R=43; C=5; V='  -1.13%  '; Force=$false},
    @{R=44; C=4; V='7.198'; Force=$true},
    @{R=44; C=5; V='  +0.61%  '; Force=$false},
    @{R=45; C=5; V='  -0.50%  '; Force=$false},
    @{R=46; C=2; V='EnergySwap'; Force=$false},
    @{R=46; C=3; V='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Force=$false},
    @{R=46; C=4; V='9.014'; Force=$true},
    @{R=46; C=5; V='  -0.46%  '; Force=$false},
    @{R=47; C=2; V='RenderToken'; Force=$false},
    @{R=47; C=3; V='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; Force=$false},
    @{R=47; C=4; V='1.687'; Force=$true},
    @{R=47; C=5; V='  +0.02%  '; Force=$false},
    @{R=48; C=2; V='Algorand'; Force=$false},
    @{R=48; C=3; V='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; Force=$false},
    @{R=48; C=4; V='0.1142'; Force=$true},
    @{R=48; C=5; V='  +1.61%  '; Force=$false},
    @{R=49; C=2; V='BabyDogeCoin'; Force=$false},
    @{R=49; C=3; V='https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; Force=$false},
    @{R=49; C=4; V='0.00000000114'; Force=$true},
    @{R=49; C=5; V='  -6.93%  '; Force=$false},
    @{R=50; C=5; V='  -0.23%  '; Force=$false},
    @{R=51; C=4; V='0.4621'; Force=$true},
    @{R=51; C=5; V='  -0.15%  '; Force=$false}
)

foreach ($e in $edits) {
    $cell = $ws.Cells.Item($e.R, $e.C)
    if ($e.Force) {
        $cell.NumberFormat = "@"
    }
    $cell.Value2 = $e.V
}
